# Auto-generated edit script: update cryptocurrency price table
# to reflect the Dec 16 2022 14:41 UTC symbol-list refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'249.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'24.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.937"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05903"
$ws.Range("D5").Style = "Normal"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "'3.425"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "5GateTokenGT"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "'6.578"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "6KuCoinTokenKCS"
$ws.Range("D8").Value = "'1.331"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.7962"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.1491"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.07783"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.03322"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.03013"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.09260"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'3.561"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.001662"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04776"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.0006026"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").Value = "'0.006207"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Value = "'0.001068"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Value = "'3.720"
$ws.Range("D23").Style = "Normal"
$ws.Range("D26").Value = "'0.1252"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'0.0006469"
$ws.Range("D27").Style = "Normal"
$ws.Range("D40").Value = "'0.04400"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.007023"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'0.1070"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'0.003364"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.01004"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.002458"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "44ACDXExchangeACXTBestin24h"
$ws.Range("D46").Value = "'0.00005890"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.9892"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'0.1103"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002098"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = "'0.01009"
$ws.Range("D51").Style = "Normal"
